# Fruta / hortaliza, semanal
# Inserts two new weekly price-report rows (172 and 173) at the top of the
# "Cebollín" data block, pushing the existing rows (old 172..218) down to
# become the new rows 174..220.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at row 172 (each insert pushes rows 172+ down by one)
$ws.Rows.Item(172).EntireRow.Insert()
$ws.Rows.Item(172).EntireRow.Insert()

# New data for row 172 (Primera) and row 173 (Segunda)
$row172 = @(7, "Terminal Hortofrutícola Agro Chillán", "Ñuble", 45215, 16, 100112037, "Cebollín", "Sin especificar", "Primera", 250, 6000, 6000, 6000, "`$/paquete 36 unidades", "Provincia de Diguillín", 167, 36, "Hortaliza")
$row173 = @(7, "Terminal Hortofrutícola Agro Chillán", "Ñuble", 45215, 16, 100112037, "Cebollín", "Sin especificar", "Segunda", 250, 5000, 5000, 5000, "`$/paquete 36 unidades", "Provincia de Diguillín", 139, 36, "Hortaliza")

for ($c = 1; $c -le 18; $c++) {
    $ws.Cells.Item(172, $c).Value2 = $row172[$c - 1]
    $ws.Cells.Item(173, $c).Value2 = $row173[$c - 1]
}
